$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each Price/Volume cell is stored as literal text (e.g. "320.65", "7.03%").
# Prefixing the assigned value with a single-quote forces Excel to keep it
# as text instead of reinterpreting it as a number/percentage, and resetting
# the cell Style back to "Normal" afterwards clears the quote-prefix marker
# so no extra formatting is left behind.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "320.65"
Set-TextValue "E2" "7.03%"
Set-TextValue "D3" "48.66"
Set-TextValue "E3" "15.37%"
Set-TextValue "D4" "5.264"
Set-TextValue "E4" "4.97%"
Set-TextValue "D5" "0.08098"
Set-TextValue "E5" "7.11%"
Set-TextValue "D6" "4.593"
Set-TextValue "E6" "4.94%"
Set-TextValue "D7" "1.644"
Set-TextValue "E7" "2.85%"
Set-TextValue "D8" "1.211"
Set-TextValue "E8" "29.05%"
Set-TextValue "D9" "0.1298"
Set-TextValue "E9" "9.18%"
Set-TextValue "E10" "5.63%"
Set-TextValue "D11" "0.09444"
Set-TextValue "E11" "4.44%"
Set-TextValue "D12" "0.04623"
Set-TextValue "E12" "10.96%"
Set-TextValue "D13" "0.1051"
Set-TextValue "E13" "0.24%"
Set-TextValue "D14" "0.001332"
Set-TextValue "E14" "3.86%"
Set-TextValue "D15" "0.04175"
Set-TextValue "E15" "1.71%"
Set-TextValue "D16" "0.005822"
Set-TextValue "E16" "0.34%"
Set-TextValue "D17" "3.338"
Set-TextValue "E17" "-0.19%"
Set-TextValue "D18" "2.430"
Set-TextValue "E18" "1.93%"
Set-TextValue "D19" "0.3419"
Set-TextValue "E19" "2.53%"
Set-TextValue "D20" "8.096"
Set-TextValue "E20" "-3.48%"
Set-TextValue "E21" "-1.26%"
Set-TextValue "D22" "0.3129"
Set-TextValue "E22" "-5.14%"
Set-TextValue "D23" "0.001308"
Set-TextValue "E23" "3.43%"
Set-TextValue "D24" "0.004247"
Set-TextValue "E24" "8.84%"
Set-TextValue "D25" "0.0001353"
Set-TextValue "E25" "6.54%"
Set-TextValue "D26" "0.0003546"
Set-TextValue "E26" "-4.78%"
Set-TextValue "D38" "0.02722"
Set-TextValue "E38" "12.98%"
Set-TextValue "D39" "0.05716"
Set-TextValue "E39" "9.32%"
Set-TextValue "D40" "0.006312"
Set-TextValue "E40" "-5.68%"
Set-TextValue "D41" "0.007770"
Set-TextValue "E41" "0.82%"
Set-TextValue "D42" "0.1442"
Set-TextValue "E42" "8.67%"
Set-TextValue "D43" "0.007718"
Set-TextValue "E43" "4.40%"
Set-TextValue "E44" "3.70%"
Set-TextValue "E45" "6.45%"
Set-TextValue "D46" "0.00006860"
Set-TextValue "E46" "9.74%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.26%"
Set-TextValue "D48" "0.06717"
Set-TextValue "E48" "46.84%"
Set-TextValue "D49" "0.004008"
Set-TextValue "E49" "-4.60%"
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.26%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.26%"
